$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5733
$ws1.Range("F3").Value = 7686
$ws1.Range("F9").Value = 4487
$ws1.Range("F11").Value = 120
$ws1.Range("F13").Value = 3027
$ws1.Range("F15").Value = 576
$ws1.Range("F16").Value = 230
$ws1.Range("F17").Value = 572
$ws1.Range("F18").Value = 490
$ws1.Range("F19").Value = 491
$ws1.Range("F22").Value = 1741
$ws1.Range("F23").Value = 1273
$ws1.Range("F25").Value = 1489
$ws1.Range("F31").Value = 31
$ws1.Range("F35").Value = 3271
$ws1.Range("F36").Value = 729
$ws1.Range("F37").Value = 54
$ws1.Range("F38").Value = 204
$ws1.Range("F40").Value = 1409

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 18

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5733
$ws4.Range("F3").Value = 7686
$ws4.Range("F9").Value = 4487
$ws4.Range("F11").Value = 120
$ws4.Range("F13").Value = 3027
$ws4.Range("F15").Value = 576
$ws4.Range("F16").Value = 230
$ws4.Range("F17").Value = 572
$ws4.Range("F18").Value = 490
$ws4.Range("F19").Value = 491
$ws4.Range("F20").Value = 18
$ws4.Range("F23").Value = 1741
$ws4.Range("F24").Value = 1273
$ws4.Range("F26").Value = 1489
$ws4.Range("F32").Value = 31
$ws4.Range("F36").Value = 3271
$ws4.Range("F38").Value = 729
$ws4.Range("F39").Value = 54
$ws4.Range("F40").Value = 204
$ws4.Range("F42").Value = 1409

$wb.Save()
